{"js": "// Remove the trailing \"Ver no Jupiter ... / \u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that used to follow the\n// \"Requisitos\" section, while leaving the \"LOM3099: Est\u00e1tica (Requisito\n// fraco)\" paragraph and everything after the footer (the blank paragraph\n// + the page-break paragraph) untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"LOM3099: Est\u00e1tica (Requisito fraco)\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// The three paragraphs immediately following the anchor are the ones the\n// commit removed: a blank \"Normal\" paragraph, the \"Ver no Jupiter...\" line,\n// and the \"\u00a9 2020 ...\" copyright line.\nconst blankParagraph = paragraphs.items[anchorIndex + 1];\nconst jupiterParagraph = paragraphs.items[anchorIndex + 2];\nconst copyrightParagraph = paragraphs.items[anchorIndex + 3];\n\nif (jupiterParagraph.text.indexOf(\"Ver no Jupiter\") === -1) {\n  throw new Error(\"Unexpected paragraph where 'Ver no Jupiter...' was expected: \" + jupiterParagraph.text);\n}\nif (copyrightParagraph.text.indexOf(\"Contact: luizeleno@usp.br\") === -1) {\n  throw new Error(\"Unexpected paragraph where the copyright line was expected: \" + copyrightParagraph.text);\n}\n\nblankParagraph.delete();\njupiterParagraph.delete();\ncopyrightParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ... / \u00a9 2020 ...\" footer block\n# (and the blank paragraph right before it) that used to follow the\n# \"Requisitos\" section, while leaving the \"LOM3099: Est\u00e1tica (Requisito\n# fraco)\" paragraph and everything after the footer (the blank paragraph\n# + the page-break paragraph) untouched.\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text\n    if ($txt -like \"LOM3099:*\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph starting with 'LOM3099:'\"\n}\n\n$blankIndex = $anchorIndex + 1\n$jupiterIndex = $anchorIndex + 2\n$copyrightIndex = $anchorIndex + 3\n\n$jupiterText = $d.Paragraphs.Item($jupiterIndex).Range.Text\nif ($jupiterText -notlike \"*Ver no Jupiter*\") {\n    throw \"Unexpected paragraph where 'Ver no Jupiter...' was expected: $jupiterText\"\n}\n$copyrightText = $d.Paragraphs.Item($copyrightIndex).Range.Text\nif ($copyrightText -notlike \"*luizeleno@usp.br*\") {\n    throw \"Unexpected paragraph where the copyright line was expected: $copyrightText\"\n}\n\n# Delete starting from the highest index so the lower indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($jupiterIndex).Range.Delete()\n$d.Paragraphs.Item($blankIndex).Range.Delete()\n"}
